# Applies the commit: "Add missing parameter entries in xlsx file and AD_IO"
# Inserts 81 new AeroDyn output-parameter rows (B1-3N1-9 Cpmin / SigCr / SgCav)
# right before the existing "B1N1AddMa" row, pushing the AddMa/Fam block (and
# everything after it) down by 81 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AeroDyn")

# Row at which the new block of 81 rows must be inserted (this is currently
# the row holding "B1N1AddMa" / "Added mass at Blade 1, Node 1").
$insertRow = 1090
$blockSize = 81

# Build the ordered list of new "Name" values: 27 Cpmin, then 27 SigCr,
# then 27 SgCav -- each following the B<blade>N<node><suffix> pattern.
$suffixes = @("Cpmin", "SigCr", "SgCav")
$names = @()
foreach ($suffix in $suffixes) {
    for ($b = 1; $b -le 3; $b++) {
        for ($n = 1; $n -le 9; $n++) {
            $names += "B" + $b + "N" + $n + $suffix
        }
    }
}

# Make room for the new rows; Excel copies the formatting of the row above
# (which already carries the correct styles for columns B and D).
$ws.Rows.Item($insertRow).Resize($blockSize).Insert()

# Fill column B (Name) first for every new row so the new "B#N#..." strings
# are appended to the shared-string table before "TODO add note" and "(-)".
for ($i = 0; $i -lt $blockSize; $i++) {
    $r = $insertRow + $i
    $ws.Cells.Item($r, 2).Value = $names[$i]
}

# Fill column D (Description) -- all new rows share the placeholder text.
for ($i = 0; $i -lt $blockSize; $i++) {
    $r = $insertRow + $i
    $ws.Cells.Item($r, 4).Value = "TODO add note"
}

# Fill column F (Units) -- all new rows use the dimensionless unit "(-)".
for ($i = 0; $i -lt $blockSize; $i++) {
    $r = $insertRow + $i
    $ws.Cells.Item($r, 6).Value = "(-)"
}

# Restore the cursor/selection to mirror the edited file's view state.
$ws.Activate()
$ws.Range("E1159").Select()
